# Auto-generated edit script applying the cell-value changes described by the diff
# against Sheets/Pandaemonium_Profits.xlsx (tabs ALC/ARM/BSM/CRP/CUL/GSM/WVR).
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H100").Value = 1790
$ws_ALC.Range("I100").Value = 1617.1428
$ws_ALC.Range("J100").Value = 3000
$ws_ALC.Range("K100").Value = 1617.1428
$ws_ALC.Range("L100").Value = 3000
$ws_ALC.Range("M100").Value = -1076.1428
$ws_ALC.Range("N100").Value = -4082
$ws_ALC.Range("H132").Value = 2209.9375
$ws_ALC.Range("I132").Value = 2209.9375
$ws_ALC.Range("K132").Value = 6629.8125
$ws_ALC.Range("M132").Value = -4099.8125
$ws_ALC.Range("H138").Value = 13341.742
$ws_ALC.Range("I138").Value = 792.2353000000001
$ws_ALC.Range("J138").Value = 28580.428
$ws_ALC.Range("K138").Value = 2376.7059
$ws_ALC.Range("L138").Value = 85741.284
$ws_ALC.Range("M138").Value = 2763.2941
$ws_ALC.Range("N138").Value = -96021.284
$ws_ALC.Range("H141").Value = 2497.4119
$ws_ALC.Range("I141").Value = 1650.9231
$ws_ALC.Range("J141").Value = 5248.5
$ws_ALC.Range("K141").Value = 4952.7693
$ws_ALC.Range("L141").Value = 15745.5
$ws_ALC.Range("M141").Value = 227.2307000000001
$ws_ALC.Range("N141").Value = -26105.5

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H24").Value = 26451.666
$ws_ARM.Range("J24").Value = 26451.666
$ws_ARM.Range("L24").Value = 26451.666
$ws_ARM.Range("N24").Value = -27199.666
$ws_ARM.Range("H37").Value = 39900
$ws_ARM.Range("I37").Value = 0
$ws_ARM.Range("K37").Value = 0
$ws_ARM.Range("M37").Value = ""
$ws_ARM.Range("H100").Value = 26451.666
$ws_ARM.Range("J100").Value = 26451.666
$ws_ARM.Range("L100").Value = 26451.666
$ws_ARM.Range("N100").Value = -28615.666
$ws_ARM.Range("H139").Value = 73171.664
$ws_ARM.Range("J139").Value = 73171.664
$ws_ARM.Range("L139").Value = 73171.664
$ws_ARM.Range("N139").Value = -83451.664

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H63").Value = 40271
$ws_BSM.Range("I63").Value = 0
$ws_BSM.Range("J63").Value = 40271
$ws_BSM.Range("K63").Value = 0
$ws_BSM.Range("L63").Value = 40271
$ws_BSM.Range("M63").Value = ""
$ws_BSM.Range("N63").Value = -41643
$ws_BSM.Range("H66").Value = 40271
$ws_BSM.Range("I66").Value = 0
$ws_BSM.Range("J66").Value = 40271
$ws_BSM.Range("K66").Value = 0
$ws_BSM.Range("L66").Value = 120813
$ws_BSM.Range("M66").Value = ""
$ws_BSM.Range("N66").Value = -127677
$ws_BSM.Range("H94").Value = 923.5714
$ws_BSM.Range("I94").Value = 1002
$ws_BSM.Range("J94").Value = 727.5
$ws_BSM.Range("K94").Value = 1002
$ws_BSM.Range("L94").Value = 727.5
$ws_BSM.Range("M94").Value = -551
$ws_BSM.Range("N94").Value = -1629.5

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H7").Value = 165.66667
$ws_CRP.Range("I7").Value = 48.5
$ws_CRP.Range("J7").Value = 400
$ws_CRP.Range("K7").Value = 48.5
$ws_CRP.Range("L7").Value = 400
$ws_CRP.Range("M7").Value = 64.5
$ws_CRP.Range("N7").Value = -626
$ws_CRP.Range("H31").Value = 8080.171
$ws_CRP.Range("I31").Value = 6948.091
$ws_CRP.Range("K31").Value = 6948.091
$ws_CRP.Range("M31").Value = -6653.091
$ws_CRP.Range("H34").Value = 8080.171
$ws_CRP.Range("I34").Value = 6948.091
$ws_CRP.Range("K34").Value = 6948.091
$ws_CRP.Range("M34").Value = -6746.091
$ws_CRP.Range("H44").Value = 8499.5
$ws_CRP.Range("I44").Value = 8333
$ws_CRP.Range("J44").Value = 8999
$ws_CRP.Range("K44").Value = 8333
$ws_CRP.Range("L44").Value = 8999
$ws_CRP.Range("M44").Value = -7891
$ws_CRP.Range("N44").Value = -9883
$ws_CRP.Range("H50").Value = 25852.572
$ws_CRP.Range("J50").Value = 25852.572
$ws_CRP.Range("L50").Value = 25852.572
$ws_CRP.Range("N50").Value = -27102.572
$ws_CRP.Range("H51").Value = 19106.666
$ws_CRP.Range("J51").Value = 19106.666
$ws_CRP.Range("L51").Value = 19106.666
$ws_CRP.Range("N51").Value = -20578.666
$ws_CRP.Range("H58").Value = 2843177
$ws_CRP.Range("I58").Value = 7577432.5
$ws_CRP.Range("J58").Value = 2623.8
$ws_CRP.Range("K58").Value = 7577432.5
$ws_CRP.Range("L58").Value = 2623.8
$ws_CRP.Range("M58").Value = -7577229.5
$ws_CRP.Range("N58").Value = -3029.8
$ws_CRP.Range("H59").Value = 29656.75
$ws_CRP.Range("I59").Value = 0
$ws_CRP.Range("J59").Value = 29656.75
$ws_CRP.Range("K59").Value = 0
$ws_CRP.Range("L59").Value = 29656.75
$ws_CRP.Range("M59").Value = ""
$ws_CRP.Range("N59").Value = -31946.75
$ws_CRP.Range("H60").Value = 8354.714
$ws_CRP.Range("J60").Value = 8672.450000000001
$ws_CRP.Range("L60").Value = 8672.450000000001
$ws_CRP.Range("N60").Value = -9694.450000000001
$ws_CRP.Range("H61").Value = 19106.666
$ws_CRP.Range("J61").Value = 19106.666
$ws_CRP.Range("L61").Value = 19106.666
$ws_CRP.Range("N61").Value = -19802.666
$ws_CRP.Range("H68").Value = 0
$ws_CRP.Range("J68").Value = 0
$ws_CRP.Range("L68").Value = 0
$ws_CRP.Range("N68").Value = ""
$ws_CRP.Range("H71").Value = 0
$ws_CRP.Range("J71").Value = 0
$ws_CRP.Range("L71").Value = 0
$ws_CRP.Range("N71").Value = ""
$ws_CRP.Range("H74").Value = 38157
$ws_CRP.Range("J74").Value = 38157
$ws_CRP.Range("L74").Value = 38157
$ws_CRP.Range("N74").Value = -39905
$ws_CRP.Range("H77").Value = 38157
$ws_CRP.Range("J77").Value = 38157
$ws_CRP.Range("L77").Value = 114471
$ws_CRP.Range("N77").Value = -123207
$ws_CRP.Range("H124").Value = 35000
$ws_CRP.Range("J124").Value = 35000
$ws_CRP.Range("L124").Value = 35000
$ws_CRP.Range("N124").Value = -39910
$ws_CRP.Range("H132").Value = 3408.25
$ws_CRP.Range("I132").Value = 2878.3333
$ws_CRP.Range("J132").Value = 4203.125
$ws_CRP.Range("K132").Value = 8634.999899999999
$ws_CRP.Range("L132").Value = 12609.375
$ws_CRP.Range("M132").Value = -6104.999899999999
$ws_CRP.Range("N132").Value = -17669.375
$ws_CRP.Range("H135").Value = 53405.715
$ws_CRP.Range("J135").Value = 53405.715
$ws_CRP.Range("L135").Value = 53405.715
$ws_CRP.Range("N135").Value = -63545.715
$ws_CRP.Range("H136").Value = 2843177
$ws_CRP.Range("I136").Value = 7577432.5
$ws_CRP.Range("J136").Value = 2623.8
$ws_CRP.Range("K136").Value = 22732297.5
$ws_CRP.Range("L136").Value = 7871.400000000001
$ws_CRP.Range("M136").Value = -22729747.5
$ws_CRP.Range("N136").Value = -12971.4

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 11912382
$ws_CUL.Range("I5").Value = 408.55554
$ws_CUL.Range("J5").Value = 33353936
$ws_CUL.Range("K5").Value = 1225.66662
$ws_CUL.Range("L5").Value = 100061808
$ws_CUL.Range("M5").Value = -1113.66662
$ws_CUL.Range("N5").Value = -100062032
$ws_CUL.Range("H70").Value = 2185.7
$ws_CUL.Range("J70").Value = 2317.4443
$ws_CUL.Range("L70").Value = 6952.3329
$ws_CUL.Range("N70").Value = -7582.3329
$ws_CUL.Range("H73").Value = 2185.7
$ws_CUL.Range("J73").Value = 2317.4443
$ws_CUL.Range("L73").Value = 6952.3329
$ws_CUL.Range("N73").Value = -9136.332900000001
$ws_CUL.Range("H107").Value = 1239.9
$ws_CUL.Range("I107").Value = 500
$ws_CUL.Range("J107").Value = 1322.1111
$ws_CUL.Range("K107").Value = 1500
$ws_CUL.Range("L107").Value = 3966.3333
$ws_CUL.Range("M107").Value = 420
$ws_CUL.Range("N107").Value = -7806.3333
$ws_CUL.Range("H122").Value = 1315.9286
$ws_CUL.Range("J122").Value = 1609
$ws_CUL.Range("L122").Value = 14481
$ws_CUL.Range("N122").Value = -19381
$ws_CUL.Range("H132").Value = 1559.8
$ws_CUL.Range("I132").Value = 1623.1538
$ws_CUL.Range("J132").Value = 1442.1428
$ws_CUL.Range("K132").Value = 14608.3842
$ws_CUL.Range("L132").Value = 12979.2852
$ws_CUL.Range("M132").Value = -12078.3842
$ws_CUL.Range("N132").Value = -18039.2852
$ws_CUL.Range("H135").Value = 11912382
$ws_CUL.Range("I135").Value = 408.55554
$ws_CUL.Range("J135").Value = 33353936
$ws_CUL.Range("K135").Value = 3676.99986
$ws_CUL.Range("L135").Value = 300185424
$ws_CUL.Range("M135").Value = -1141.99986
$ws_CUL.Range("N135").Value = -300190494

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H80").Value = 43000
$ws_GSM.Range("J80").Value = 43000
$ws_GSM.Range("L80").Value = 43000
$ws_GSM.Range("N80").Value = -44996
$ws_GSM.Range("H83").Value = 43000
$ws_GSM.Range("J83").Value = 43000
$ws_GSM.Range("L83").Value = 215000
$ws_GSM.Range("N83").Value = -224984
$ws_GSM.Range("H92").Value = 0
$ws_GSM.Range("J92").Value = 0
$ws_GSM.Range("L92").Value = 0
$ws_GSM.Range("N92").Value = ""
$ws_GSM.Range("H97").Value = 0
$ws_GSM.Range("I97").Value = 0
$ws_GSM.Range("J97").Value = 0
$ws_GSM.Range("K97").Value = 0
$ws_GSM.Range("L97").Value = 0
$ws_GSM.Range("M97").Value = ""
$ws_GSM.Range("N97").Value = ""
$ws_GSM.Range("H141").Value = 40854.375
$ws_GSM.Range("J141").Value = 40854.375
$ws_GSM.Range("L141").Value = 40854.375
$ws_GSM.Range("N141").Value = -51214.375

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H138").Value = 44424.777
$ws_WVR.Range("J138").Value = 44424.777
$ws_WVR.Range("L138").Value = 44424.777
$ws_WVR.Range("N138").Value = -54704.777
$ws_WVR.Range("H140").Value = 48593.855
$ws_WVR.Range("J140").Value = 48593.855
$ws_WVR.Range("L140").Value = 48593.855
$ws_WVR.Range("N140").Value = -58953.855
$ws_WVR.Range("H141").Value = 44900
$ws_WVR.Range("J141").Value = 44900
$ws_WVR.Range("L141").Value = 44900
$ws_WVR.Range("N141").Value = -55260
